$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on Price (D) and Volume (E) columns so numeric-looking
# strings (e.g. "1.00", "129.89") are stored as text, matching the source data.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "59.022.33"
$ws.Range("E2").Value = "  -0.04%  "
$ws.Range("D3").Value = "2.511.46"
$ws.Range("E3").Value = "  +0.20%  "
$ws.Range("E4").Value = "  +0.18%  "
$ws.Range("D5").Value = "533.23"
$ws.Range("E5").Value = "  -0.47%  "
$ws.Range("D6").Value = "135.75"
$ws.Range("E6").Value = "  -1.07%  "
$ws.Range("D7").Value = "1.00"
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("E8").Value = "  +0.11%  "
$ws.Range("D9").Value = "0.101"
$ws.Range("E9").Value = "  +0.32%  "
$ws.Range("D10").Value = "0.158"
$ws.Range("E10").Value = "  -1.28%  "
$ws.Range("E11").Value = "  +1.29%  "
$ws.Range("D12").Value = "0.345"
$ws.Range("E12").Value = "  -0.60%  "
$ws.Range("D13").Value = "2.957.99"
$ws.Range("E13").Value = "  -0.36%  "
$ws.Range("D14").Value = "58.918.32"
$ws.Range("E14").Value = "  -0.33%  "
$ws.Range("D15").Value = "22.81"
$ws.Range("E15").Value = "  -1.76%  "
$ws.Range("E16").Value = "  -1.09%  "
$ws.Range("D17").Value = "2.513.62"
$ws.Range("E17").Value = "  +0.11%  "
$ws.Range("D18").Value = "11.04"
$ws.Range("E18").Value = "  -0.27%  "
$ws.Range("D19").Value = "4.25"
$ws.Range("E19").Value = "  -0.39%  "
$ws.Range("D20").Value = "323.77"
$ws.Range("E20").Value = "  -0.52%  "
$ws.Range("E21").Value = "  -0.17%  "
$ws.Range("D22").Value = "5.92"
$ws.Range("E22").Value = "  +0.99%  "
$ws.Range("D23").Value = "65.20"
$ws.Range("E23").Value = "  +0.18%  "
$ws.Range("E24").Value = "  -0.15%  "
$ws.Range("E25").Value = "  -1.03%  "
$ws.Range("D26").Value = "0.998"
$ws.Range("E26").Value = "  -0.27%  "
$ws.Range("E27").Value = "  -0.81%  "
$ws.Range("D28").Value = "0.0₃0764"
$ws.Range("E28").Value = "  -1.38%  "
$ws.Range("D29").Value = "6.48"
$ws.Range("E29").Value = "  -3.74%  "
$ws.Range("E30").Value = "  -1.32%  "
$ws.Range("D31").Value = "169.22"
$ws.Range("E31").Value = "  +1.14%  "
$ws.Range("E32").Value = "  +0.00%  "
$ws.Range("E33").Value = "  -3.94%  "
$ws.Range("B34").Value = "ImmutableX"
$ws.Range("C34").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D34").Value = "1.37"
$ws.Range("E34").Value = "  -2.63%  "
$ws.Range("B35").Value = "EthereumClassic"
$ws.Range("C35").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D35").Value = "18.40"
$ws.Range("E35").Value = "  -0.95%  "
$ws.Range("D36").Value = "4.05"
$ws.Range("E36").Value = "  -1.78%  "
$ws.Range("E37").Value = "  -2.92%  "
$ws.Range("D38").Value = "0.805"
$ws.Range("E38").Value = "  -3.20%  "
$ws.Range("D39").Value = "3.58"
$ws.Range("E39").Value = "  -1.63%  "
$ws.Range("D40").Value = "282.45"
$ws.Range("E40").Value = "  +0.70%  "
$ws.Range("E41").Value = "  +0.18%  "
$ws.Range("E42").Value = "  -4.28%  "
$ws.Range("B43").Value = "Mantle"
$ws.Range("C43").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D43").Value = "0.604"
$ws.Range("E43").Value = "  -0.07%  "
$ws.Range("B44").Value = "Aave"
$ws.Range("C44").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D44").Value = "129.89"
$ws.Range("E44").Value = "  +1.52%  "
$ws.Range("E45").Value = "  +0.17%  "
$ws.Range("E46").Value = "  -0.45%  "
$ws.Range("D47").Value = "0.0501"
$ws.Range("E47").Value = "  -2.31%  "
$ws.Range("D48").Value = "0.0218"
$ws.Range("E48").Value = "  -2.06%  "
$ws.Range("D49").Value = "17.30"
$ws.Range("E49").Value = "  -0.41%  "
$ws.Range("D50").Value = "1.758.48"
$ws.Range("E50").Value = "  -0.70%  "
$ws.Range("E51").Value = "  -0.51%  "

# Restore default cell style (the text format above would otherwise leave
# the cells tagged with a distinct style index).
$ws.Range("D2:E51").Style = "Normal"
